$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit cyclically rotates the record data held in rows 33-36:
#   row 36's record moves up to row 33
#   row 33's record moves down to row 34
#   row 34's record moves down to row 35
#   row 35's record moves down to row 36
# Only the columns whose values actually differ between these four rows are
# touched (id/taxon names/author/coords/observers + the optional
# "Enhet/Ålder-Stadium/Kön/Aktivitet/Metod/Bestämningsmetod" block), so
# everything else (location text, dates, flags, etc. - identical across the
# four rows) is left completely untouched.

$firstRow = 33
$lastRow = 36

# Columns (1-based) that carry data which differs row-to-row and therefore
# needs to move with the rotation.
# A=1 B=2 D=4 E=5 F=6 G=7 H=8 J=10 K=11 L=12 M=13 N=14 Q=17 R=18 AF=32 AX=50
$cols = @(1, 2, 4, 5, 6, 7, 8, 10, 11, 12, 13, 14, 17, 18, 32, 50)

# Snapshot current values AND presence (present-but-blank vs not-present-at-all
# are distinct on this engine: a cell with no value at all reads back as
# PowerShell $null, while an existing-but-empty cell reads back as "").
$snapshot = @{}
foreach ($r in $firstRow..$lastRow) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# destination row -> source row (rotate down by one, wrapping at the bottom)
$mapping = @{
    33 = 36
    34 = 33
    35 = 34
    36 = 35
}

foreach ($destRow in 33..36) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    $curVals = $snapshot[$destRow]

    foreach ($c in $cols) {
        $newVal = $srcVals[$c]
        $curVal = $curVals[$c]

        $curPresent = ($curVal -ne $null)
        $newPresent = ($newVal -ne $null)

        if ($newPresent -eq $curPresent -and $newVal -eq $curVal) {
            # Nothing to do - writing "" to an already-blank-and-present cell
            # (or re-writing an identical value) is a no-op we must skip,
            # because on this engine assigning "" to an existing cell clears
            # it away entirely (same as assigning $null), which would wrongly
            # turn an already-correct present-but-blank cell into an absent
            # one.
            continue
        }

        if (-not $newPresent) {
            # Cell must become entirely absent.
            $ws.Cells.Item($destRow, $c).Value2 = $null
        }
        else {
            # Cell must exist (blank or with real content) - a plain
            # assignment both creates missing cells and updates existing
            # non-blank ones without side effects.
            $ws.Cells.Item($destRow, $c).Value2 = $newVal
        }
    }
}
